# Generate Report for Handback
# Updates the "Latest Handback DateTime" column (column G) on each of the
# per-language worksheets. For each sheet, every cell in column G that
# currently holds the sheet's old handback timestamp text is updated to the
# new timestamp text produced by the (re-)generated handback report.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Sheet = "zh-cn"; OldValue = "2016-02-22 08:48:32"; NewValue = "2016-02-22 08:49:28" },
    @{ Sheet = "de-de"; OldValue = "2016-02-22 08:48:45"; NewValue = "2016-02-22 08:49:38" },
    @{ Sheet = "ja-jp"; OldValue = "2016-02-22 08:48:57"; NewValue = "2016-02-22 08:49:48" },
    @{ Sheet = "zh-tw"; OldValue = "2016-02-22 08:49:07"; NewValue = "2016-02-22 08:49:58" }
)

foreach ($update in $updates) {
    $ws = $wb.Sheets.Item($update.Sheet)
    $usedRange = $ws.UsedRange
    $rowCount = $usedRange.Rows.Count

    for ($r = 1; $r -le $rowCount; $r++) {
        $cell = $ws.Cells.Item($r, 7)  # Column G = "Latest Handback DateTime"
        if ($cell.Text -eq $update.OldValue) {
            $cell.Value = $update.NewValue
        }
    }
}
